$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the two publish-date header cells (shared strings)
$ws.Range("I9").Value = "1402-04-06 (11)"
$ws.Range("M9").Value = "1402-04-06 (3)"

# Update the cumulative figures in column M (12-month period ended 1401/12)
$ws.Range("M14").Value = -27471
$ws.Range("M17").Value = 35823
$ws.Range("M19").Value = 77462
$ws.Range("M20").Value = 91124
$ws.Range("M21").Value = -2001
$ws.Range("M22").Value = 89123
$ws.Range("M24").Value = 89123
